# Add columns I (I0) and J (IF) to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - style matches the existing header row (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-24
$data = @(
    @(8, 9),   # row 2
    @(8, 8),   # row 3
    @(9, 9),   # row 4
    @(8, 9),   # row 5
    @(8, 8),   # row 6
    @(9, 9),   # row 7
    @(7, 8),   # row 8
    @(7, 8),   # row 9
    @(1, 2),   # row 10
    @(8, 8),   # row 11
    @(7, 8),   # row 12
    @(9, 9),   # row 13
    @(1, 1),   # row 14
    @(9, 9),   # row 15
    @(6, 6),   # row 16
    @(6, 6),   # row 17
    @(7, 8),   # row 18
    @(6, 6),   # row 19
    @(8, 8),   # row 20
    @(6, 6),   # row 21
    @(7, 7),   # row 22
    @(3, 3),   # row 23
    @(2, 2)    # row 24
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 9).Value = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $r++
}
